$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "uliza"
$ws.Range("B2").Value = "ZAR 7000"

# Force the date cell to be stored as plain text, not an Excel date serial
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-01-30"

$ws.Range("D2").Value = "Pending"

# Remove row 3 (was BMW / ZAR 10000) entirely
$ws.Range("A3:B3").Delete()
